$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.886.89"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.550.94"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'206.08"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "'21.50"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.772.53"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.558.76"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "26.901.61"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'61.65"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "'213.47"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").Value = "'153.05"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'6.66"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "1.373.33"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "'0.968"
$ws.Range("E36").Value = "  +6.30%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").Value = "'0.522"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("D43").Value = "'5.49"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +3.56%  "
$ws.Range("D45").Value = "'63.59"
$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  -2.29%  "
$ws.Range("D47").Value = "1.686.45"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "'86.22"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "'0.0508"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  +0.38%  "
